# Apply edits described by the commit diff:
# 1. Window width of the workbook view changes from 18345 to 27945.
# 2. Shared string "Petit Lotu" (Secteur column, rows 10-13) renamed to "Mezzanu".
# 3. Shared string "Trave" (Secteur column, rows 26-29) renamed to "U Travu".
# 4. The active selection on the sheet moves from K10 to B26:B29 (activeCell B26).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meteo")

# --- 2 & 3: update the "Secteur" values that changed name ---
$ws.Range("B10:B13").Value = "Mezzanu"
$ws.Range("B26:B29").Value = "U Travu"

# --- 4: move / set the active selection to B26:B29 with active cell B26 ---
$ws.Activate()
$ws.Range("B26:B29").Select()

# --- 1: widen the workbook window ---
$excel.ActiveWindow.Width = 27945 / 20
$excel.ActiveWindow.Height = 12180 / 20
